$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (ano 2021)
$ws.Range("C2").Value = 484
$ws.Range("D2").Value = 65
$ws.Range("E2").Value = 419
$ws.Range("F2").Value = 60.74766355140186
$ws.Range("G2").Value = 86.57024793388429
$ws.Range("H2").Value = 13.4297520661157

# Row 3 (ano 2022)
$ws.Range("C3").Value = 362
$ws.Range("D3").Value = 167
$ws.Range("E3").Value = 195
$ws.Range("F3").Value = 34.50413223140496
$ws.Range("G3").Value = 53.86740331491713
$ws.Range("H3").Value = 46.13259668508287

# Row 4 (ano 2023)
$ws.Range("C4").Value = 320
$ws.Range("D4").Value = 186
$ws.Range("E4").Value = 134
$ws.Range("F4").Value = 51.38121546961326
$ws.Range("G4").Value = 41.875
$ws.Range("H4").Value = 58.12500000000001

# Row 5 (ano 2024)
$ws.Range("C5").Value = 457
$ws.Range("D5").Value = 233
$ws.Range("E5").Value = 224
$ws.Range("F5").Value = 72.8125
$ws.Range("G5").Value = 49.01531728665208
$ws.Range("H5").Value = 50.98468271334792

# Row 6 (ano 2025)
$ws.Range("C6").Value = 386
$ws.Range("D6").Value = 305
$ws.Range("E6").Value = 81
$ws.Range("F6").Value = 66.73960612691467
$ws.Range("G6").Value = 20.98445595854922
$ws.Range("H6").Value = 79.01554404145078
